{"js": "// Convert the two \"simple\" field codes (w:fldSimple) that mark the user\n// content zone (m:usercontent zone1 / m:endusercontent) into the\n// \"complex\" field representation (begin/instrText/separate/end field\n// characters spread across runs). This mirrors what Word itself does the\n// first time it saves a document containing fldSimple fields, and is the\n// fix for \"user content is lost after two generations without edition\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Replaces the lone w:fldSimple field inside `paragraph` with the\n// equivalent complex field (4 runs: begin / instrText / separate / end),\n// leaving the paragraph element itself (and its properties) untouched.\nasync function convertSimpleFieldToComplex(paragraph, instrText) {\n    const startRange = paragraph.getRange(\"Start\");\n\n    const ooxml =\n        '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p>' +\n        '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n        '<w:r><w:instrText>' + instrText + '</w:instrText></w:r>' +\n        '<w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n        '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part>' +\n        '</pkg:package>';\n\n    // Insert the new complex-field runs right before the existing content\n    // of the paragraph (this keeps the original <w:p> element - and its\n    // rsid/style attributes - untouched).\n    startRange.insertOoxml(ooxml, Word.InsertLocation.before);\n    await context.sync();\n\n    // The paragraph now contains two fields with the same instruction\n    // text: the newly inserted complex one and the original fldSimple.\n    // Remove the original (now last) one.\n    const fields = paragraph.fields;\n    fields.load(\"items\");\n    await context.sync();\n\n    fields.items[fields.items.length - 1].delete();\n    await context.sync();\n}\n\n// Paragraph index 1: \"m:usercontent zone1\"\nconst userContentParagraph = paragraphs.items[1];\nawait convertSimpleFieldToComplex(userContentParagraph, \"m:usercontent zone1\");\n\n// Paragraph index 3: \"m:endusercontent\"\nconst endUserContentParagraph = paragraphs.items[3];\nawait convertSimpleFieldToComplex(endUserContentParagraph, \"m:endusercontent\");\n", "ps1": "# Convert the two \"simple\" field codes (w:fldSimple) that delimit the\n# user content zone (m:usercontent zone1 / m:endusercontent) into the\n# \"complex\" field representation (begin/instrText/separate/end field\n# characters spread across runs). This mirrors what Word itself does the\n# first time it saves a document containing fldSimple fields, and is the\n# fix for \"user content is lost after two generations without edition\".\n\n$d = $word.ActiveDocument\n\n# Range.Paragraphs / Range.Fields on a zero-length or otherwise scoped\n# sub-range are not reliable in this engine (they can report the whole\n# -document paragraph/field count and hand back an unrelated item), so\n# paragraphs are located by scanning Document.Paragraphs and comparing\n# Start/End bounds, and fields are always addressed through\n# Document.Fields using their 1-based document-order Index.\nfunction Get-ParagraphAtPosition($pos) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {\n            return $p\n        }\n    }\n    return $d.Paragraphs.Item($d.Paragraphs.Count)\n}\n\n# Replaces the field at `$fieldIndex` (1-based index into Document.Fields)\n# with the equivalent complex field (4 runs: begin / instrText / separate\n# / end) inserted at the start of its paragraph, then removes the\n# original fldSimple - leaving the paragraph element itself (and its\n# properties) untouched.\nfunction Convert-FieldAtIndexToComplex($fieldIndex) {\n    $fieldObj = $d.Fields.Item($fieldIndex)\n    $instr = $fieldObj.Code.Text.Trim()\n\n    $codeStart = $fieldObj.Code.Start\n    $paraRange = (Get-ParagraphAtPosition $codeStart).Range\n    $startPos = $paraRange.Start\n    $insertionRange = $d.Range($startPos, $startPos)\n\n    $ooxmlFrag = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n        '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' +\n        '<w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n        '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n        '</w:p>'\n\n    $insertionRange.InsertXML($ooxmlFrag)\n\n    # The newly inserted complex field takes the original field's slot in\n    # document order; the original fldSimple is pushed one slot later.\n    # Delete it there.\n    $d.Fields.Item($fieldIndex + 1).Delete()\n}\n\n# Identify target field indices (1-based, Document.Fields order) up\n# front, before any mutation.\n$targetIndices = New-Object System.Collections.ArrayList\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n    $f = $d.Fields.Item($i)\n    $code = $f.Code.Text.Trim()\n    if ($code -eq \"m:usercontent zone1\" -or $code -eq \"m:endusercontent\") {\n        [void]$targetIndices.Add($i)\n    }\n}\n\n# Process from the last field to the first: converting a field only\n# inserts/removes content at or after its own Document.Fields index, so\n# walking back-to-front keeps the not-yet-processed (earlier) indices\n# valid.\nfor ($t = $targetIndices.Count - 1; $t -ge 0; $t--) {\n    Convert-FieldAtIndexToComplex $targetIndices[$t]\n}\n"}
